$d = $word.ActiveDocument

# The document ends with the paragraph "...material has no meaning." which
# is immediately followed by the _GoBack bookmark and then the section
# properties. We append a new "Event System" section after it: a Heading 2
# title, a body paragraph, a blank separator paragraph, and a final body
# paragraph.

# Create all four new paragraphs first, chaining off the still-plain
# (non-heading) ranges so the Heading 2 style we apply afterwards doesn't
# "leak" forward into the paragraphs that follow it.
$anchor = $d.Paragraphs.Last
$r0 = $anchor.Range
$r0.Collapse(0)
$r0.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last

$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last

$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last

$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last

# Paragraph 1: new "Event System" Heading 2 title
$p1.Range.Text = "Event System"
$p1.Style = "Heading 2"

# Paragraph 2: body text about key input constants
$p2.Range.Text = "Define constants for individual key inputs. For each type of context API, have an array that maps the API’s keycodes to the universal key inputs. Do similar for states."

# Paragraph 3: stays a blank separator line. InsertParagraphAfter leaves a
# leftover empty run at the split point, so type a placeholder and remove
# it again with an explicit Range.Delete to end up with a clean, run-less
# empty paragraph (matching how blank paragraphs look elsewhere in this
# document).
$p3.Range.Text = "X"
$p3Start = $p3.Range.Start
$placeholder = $d.Range($p3Start, $p3Start + 1)
$placeholder.Delete()

# Paragraph 4: body text about event interfaces
$p4.Range.Text = "Components can implement event interfaces, when they do, the interface must subscribe to the event object."

Write-Host "Added Event System section; document now has" $d.Paragraphs.Count "paragraphs"
